$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Rows("6:9").Insert()
$ws.Rows("29:35").Insert()

$ws.Range("A6").Value = "However, the largest currently operating DAC facilities use solid sorbent systems and utilize non-fossil resources at low heat (80-120 C),"
$ws.Range("A7").Value = "see: https://www.wri.org/insights/direct-air-capture-resource-considerations-and-costs-carbon-removal, Status of the Leading DAC Companies."
$ws.Range("A8").Value = "Therefore, we take the heat input specified for DAC1 plants and convert to the amount of electricity needed"
$ws.Range("A9").Value = "if supplied by a heat pump."

